{"js": "// Update the Introduction paragraph of the Big Mountain Resort report:\n//   1. Insert a new sentence about the $1.5M additional annual chairlift\n//      cost right after \"...market segment. \" (and lower-case the \"The\"\n//      that begins the following sentence, since it's now mid-sentence).\n//   2. Remove \"by at least 15% \" from \"...increase the revenue by at\n//      least 15% at the end...\".\n\nconst body = context.document.body;\n\n// --- Step 1: insert the new sentence after \"market segment. \" ---\nconst anchor1 = body.search(\"market segment. The Big Mountain is looking\", {\n  matchCase: true\n});\nanchor1.load(\"text\");\nawait context.sync();\n\nif (anchor1.items.length > 0) {\n  const target1 = anchor1.items[0];\n  const replacement1 =\n    \"market segment. There is a additional $1.5M annual cost after the Big \" +\n    \"Mountain newly installed a chairlift. Therefore, the Big Mountain is looking\";\n  target1.insertText(replacement1, \"Replace\");\n  await context.sync();\n}\n\n// --- Step 2: remove \"by at least 15% \" ---\nconst anchor2 = body.search(\"revenue by at least 15% at the end\", {\n  matchCase: true\n});\nanchor2.load(\"text\");\nawait context.sync();\n\nif (anchor2.items.length > 0) {\n  const target2 = anchor2.items[0];\n  target2.insertText(\"revenue at the end\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Update the Introduction paragraph of the Big Mountain Resort report:\n#   1. Insert a new sentence about the $1.5M additional annual chairlift\n#      cost right after \"...market segment. \" (and lower-case the \"The\"\n#      that begins the following sentence, since it's now mid-sentence).\n#   2. Remove \"by at least 15% \" from \"...increase the revenue by at\n#      least 15% at the end...\".\n\n$d = $word.ActiveDocument\n\n# --- Step 1: insert the new sentence after \"market segment. \" ---\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"market segment. The Big Mountain is looking\"\n$find1.Replacement.Text = \"market segment. There is a additional `$1.5M annual cost after the Big Mountain newly installed a chairlift. Therefore, the Big Mountain is looking\"\n$find1.Forward = $true\n$find1.Wrap = 1\n$result1 = $find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)\n\n# --- Step 2: remove \"by at least 15% \" ---\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"revenue by at least 15% at the end\"\n$find2.Replacement.Text = \"revenue at the end\"\n$find2.Forward = $true\n$find2.Wrap = 1\n$result2 = $find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n\nWrite-Output \"step1=$result1 step2=$result2\"\n"}
